# add notebook part1 cleaning data
$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsClean = $wb.Worksheets.Item("Donnée propre")
$wsClean.Name = "0. Data cleaned"

$wsDirty = $wb.Worksheets.Item("Donnée à nettoyer")
$wsDirty.Name = "1. Data to clean"

# --- Fix a couple of mis-typed values on the "to clean" sheet ---
# Row 2: country should read "USA" (not "USA_")
$wsDirty.Range("C2").Value = "USA"

# Row 4: currency label should read "euros" (lowercase) here
$wsDirty.Range("G4").Value = "euros"

# Row 17: currency label should read "EUR" here
$wsDirty.Range("G17").Value = "EUR"

# --- Add the new "contact" column (I) with a couple of phone numbers ---
$wsDirty.Range("I1").Value = "contact"
$wsDirty.Columns.Item(9).ColumnWidth = 13.33

$wsDirty.Range("I9").Value = 5960023
$wsDirty.Range("I17").Value = 35206512

# --- Update selections to match the edited workbook ---
$wsClean.Range("I9").Select()
$wsDirty.Range("H9").Select()

$wsDirty.Activate()
